$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.105.84'
$ws.Range("E2").Value = '  -3.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.955.24'
$ws.Range("E3").Value = '  -3.36%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.19'
$ws.Range("E5").Value = '  -12.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.593'
$ws.Range("E6").Value = '  -4.27%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.42'
$ws.Range("E8").Value = '  -6.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.363'
$ws.Range("E9").Value = '  -5.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.78'
$ws.Range("E10").Value = '  +0.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0738'
$ws.Range("E11").Value = '  -6.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0976'
$ws.Range("E12").Value = '  -4.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.240.94'
$ws.Range("E13").Value = '  -3.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.67'
$ws.Range("E14").Value = '  -6.43%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.49'
$ws.Range("E15").Value = '  -6.36%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.742'
$ws.Range("E16").Value = '  -8.83%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.93'
$ws.Range("E17").Value = '  -7.22%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.954.61'
$ws.Range("E18").Value = '  -4.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.065.01'
$ws.Range("E19").Value = '  -3.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.09'
$ws.Range("E20").Value = '  -4.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0789'
$ws.Range("E21").Value = '  -7.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.96'
$ws.Range("E22").Value = '  -5.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '220.73'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.32'
$ws.Range("E26").Value = '  -13.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.37'
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.38'
$ws.Range("E28").Value = '  -7.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.55'
$ws.Range("E29").Value = '  -6.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").Value = '  -3.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.121'
$ws.Range("E31").Value = '  -7.63%  '

$ws.Range("E32").Value = '  -4.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.29'
$ws.Range("E33").Value = '  -8.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0595'
$ws.Range("E34").Value = '  -10.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  -9.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.28'
$ws.Range("E36").Value = '  -6.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -2.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("E39").Value = '  -7.70%  '

$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.13'
$ws.Range("E40").Value = '  -2.10%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.98'
$ws.Range("E41").Value = '  -2.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.401.10'
$ws.Range("E42").Value = '  -1.24%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0867'
$ws.Range("E43").Value = '  -9.70%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0197'
$ws.Range("E44").Value = '  -8.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.07'
$ws.Range("E45").Value = '  -11.93%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.00'
$ws.Range("E46").Value = '  -5.73%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '14.53'
$ws.Range("E47").Value = '  -9.37%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.970'
$ws.Range("E48").Value = '  -6.26%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.84'
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("B50").Value = 'FTXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.65'
$ws.Range("E50").Value = '  +16.82%  '

$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.61'
$ws.Range("E51").Value = '  -10.33%  '
